$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-15
# from serial date 45204 (2023-10-05) to 45207 (2023-10-08)
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
